$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Rename the sheet from "Data" to "Summary"
# ------------------------------------------------------------------
$ws.Name = "Summary"

# ------------------------------------------------------------------
# 1b. Re-assert the formatting of the two untouched header cells so
#     their font attributes are preserved through the save (the
#     "Afghanistan" title stays 18pt, the section heading stays bold).
# ------------------------------------------------------------------
$ws.Range("A1").Font.Size = 18
$ws.Range("A1").Font.Bold = $false

$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Size = 11

# ------------------------------------------------------------------
# 1c. Add a new named cell style "title_" (bold + underline) - this is
#     the new style introduced between "title" and "source" in the
#     workbook's style sheet.
# ------------------------------------------------------------------
$titleStyle = $wb.Styles.Add("title_")
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true

# ------------------------------------------------------------------
# 2. Remove the old rows 5-7 (their content is being replaced by a
#    new, more detailed layout further down the sheet). Delete from
#    the bottom up so row numbers of not-yet-deleted rows stay valid.
# ------------------------------------------------------------------
$ws.Range("A7:D7").EntireRow.Delete()
$ws.Range("A6:D6").EntireRow.Delete()
$ws.Range("A5:D5").EntireRow.Delete()

# ------------------------------------------------------------------
# 3. Rebuild the sheet content in its new layout.
#    Rows 1 and 3 (Afghanistan / MSME Participation on the Economy)
#    are untouched.
# ------------------------------------------------------------------

# Row 9: new bold+underlined sub-heading (uses the new "title_" style)
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Style = "title_"

# Row 11: column headers (bold), same style previously used on row 5
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true

$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true

$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# Row 12: new data row - Employment (% of total) / 50
$ws.Range("A12").Value = "Employment (% of total)"
$ws.Range("A12").Font.Bold = $true

$ws.Range("D12").Value = "'50"

# Row 13: Enterprises (% of total) / 85 (moved down from old row 6)
$ws.Range("A13").Value = "Enterprises (% of total)"
$ws.Range("A13").Font.Bold = $true

$ws.Range("D13").Value = "'85"

# Row 14: source note (italic), moved down from old row 7
$ws.Range("A14").Value = "Source: MFA, 2010"
$ws.Range("A14").Font.Italic = $true

# Row 20: bold "MFA" label
$ws.Range("A20").Value = "MFA"
$ws.Range("A20").Font.Bold = $true

# Row 21: full citation text (italic)
$ws.Range("A21").Value = "Ministry of Foreign Affaris (MFA), ""SMALL TO MEDIUM ENTERPRISE PAPERS"", N/S, p. 3. Available at http://mfa.gov.af/content/files/SME%20PAPER.pdf"
$ws.Range("A21").Font.Italic = $true
